$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Cells.Item(3, 7).Value = 2.9  # G3: 2.7 -> 2.9
$ws.Cells.Item(3, 8).Value = 3.5  # H3: 3.4 -> 3.5
$ws.Cells.Item(3, 9).Value = 2.38  # I3: 2.6 -> 2.38
$ws.Cells.Item(3, 10).Value = 3.2  # J3: 3.1 -> 3.2
$ws.Cells.Item(3, 12).Value = 2.88  # L3: 3 -> 2.88
$ws.Cells.Item(3, 25).Value = 11  # Y3: 10 -> 11
$ws.Cells.Item(3, 26).Value = 29  # Z3: 26 -> 29
$ws.Cells.Item(3, 27).Value = 21  # AA3: 19 -> 21
$ws.Cells.Item(3, 32).Value = 34  # AF3: 29 -> 34
$ws.Cells.Item(3, 34).Value = 12  # AH3: 13 -> 12
$ws.Cells.Item(3, 36).Value = 9.5  # AJ3: 10 -> 9.5
$ws.Cells.Item(3, 37).Value = 23  # AK3: 26 -> 23
$ws.Cells.Item(3, 38).Value = 17  # AL3: 19 -> 17
$ws.Cells.Item(3, 41).Value = 15  # AO3: 13 -> 15
$ws.Cells.Item(3, 51).Value = 12  # AY3: 13 -> 12

# Row 5
$ws.Cells.Item(5, 7).Value = 2.25  # G5: 2.5 -> 2.25
$ws.Cells.Item(5, 9).Value = 3.4  # I5: 2.9 -> 3.4
$ws.Cells.Item(5, 10).Value = 3.1  # J5: 3.4 -> 3.1
$ws.Cells.Item(5, 12).Value = 4.33  # L5: 3.75 -> 4.33
$ws.Cells.Item(5, 15).Value = 1.5  # O5: 1.44 -> 1.5
$ws.Cells.Item(5, 16).Value = 2.5  # P5: 2.63 -> 2.5
$ws.Cells.Item(5, 19).Value = 1.57  # S5: 1.53 -> 1.57
$ws.Cells.Item(5, 20).Value = 2.25  # T5: 2.38 -> 2.25
$ws.Cells.Item(5, 23).Value = 6  # W5: 6.5 -> 6
$ws.Cells.Item(5, 24).Value = 9.5  # X5: 11 -> 9.5
$ws.Cells.Item(5, 25).Value = 10  # Y5: 11 -> 10
$ws.Cells.Item(5, 26).Value = 21  # Z5: 26 -> 21
$ws.Cells.Item(5, 27).Value = 21  # AA5: 23 -> 21
$ws.Cells.Item(5, 31).Value = 19  # AE5: 17 -> 19
$ws.Cells.Item(5, 34).Value = 8  # AH5: 7.5 -> 8
$ws.Cells.Item(5, 35).Value = 15  # AI5: 13 -> 15
$ws.Cells.Item(5, 36).Value = 13  # AJ5: 12 -> 13
$ws.Cells.Item(5, 37).Value = 41  # AK5: 29 -> 41
$ws.Cells.Item(5, 38).Value = 34  # AL5: 29 -> 34
$ws.Cells.Item(5, 40).Value = 4  # AN5: 4.33 -> 4
$ws.Cells.Item(5, 41).Value = 13  # AO5: 15 -> 13
$ws.Cells.Item(5, 46).Value = 2.25  # AT5: 2.38 -> 2.25
$ws.Cells.Item(5, 50).Value = 5  # AX5: 4.75 -> 5
$ws.Cells.Item(5, 51).Value = 21  # AY5: 19 -> 21
$ws.Cells.Item(5, 53).Value = 67  # BA5: 51 -> 67
$ws.Cells.Item(5, 55).Value = 351  # BC5: 301 -> 351

# Row 6
$ws.Cells.Item(6, 7).Value = 3.6  # G6: 4 -> 3.6
$ws.Cells.Item(6, 8).Value = 3.1  # H6: 3.2 -> 3.1
$ws.Cells.Item(6, 9).Value = 2.1  # I6: 1.95 -> 2.1
$ws.Cells.Item(6, 10).Value = 4.33  # J6: 4.75 -> 4.33
$ws.Cells.Item(6, 12).Value = 2.88  # L6: 2.75 -> 2.88
$ws.Cells.Item(6, 13).Value = 1.1  # M6: 1.08 -> 1.1
$ws.Cells.Item(6, 14).Value = 7  # N6: 8 -> 7
$ws.Cells.Item(6, 17).Value = 2.5  # Q6: 2.4 -> 2.5
$ws.Cells.Item(6, 18).Value = 1.5  # R6: 1.53 -> 1.5
$ws.Cells.Item(6, 21).Value = 2.1  # U6: 2.2 -> 2.1
$ws.Cells.Item(6, 22).Value = 1.67  # V6: 1.62 -> 1.67
$ws.Cells.Item(6, 23).Value = 8.5  # W6: 9 -> 8.5
$ws.Cells.Item(6, 24).Value = 17  # X6: 19 -> 17
$ws.Cells.Item(6, 25).Value = 13  # Y6: 15 -> 13
$ws.Cells.Item(6, 27).Value = 34  # AA6: 41 -> 34
$ws.Cells.Item(6, 28).Value = 41  # AB6: 51 -> 41
$ws.Cells.Item(6, 30).Value = 6  # AD6: 6.5 -> 6
$ws.Cells.Item(6, 34).Value = 6  # AH6: 5.5 -> 6
$ws.Cells.Item(6, 35).Value = 9  # AI6: 8 -> 9
$ws.Cells.Item(6, 37).Value = 19  # AK6: 17 -> 19
$ws.Cells.Item(6, 38).Value = 21  # AL6: 19 -> 21
$ws.Cells.Item(6, 40).Value = 5.5  # AN6: 6 -> 5.5
$ws.Cells.Item(6, 41).Value = 21  # AO6: 23 -> 21
$ws.Cells.Item(6, 42).Value = 34  # AP6: 41 -> 34
$ws.Cells.Item(6, 45).Value = 301  # AS6: 351 -> 301
$ws.Cells.Item(6, 50).Value = 4  # AX6: 3.75 -> 4
$ws.Cells.Item(6, 51).Value = 12  # AY6: 11 -> 12

# Row 7
$ws.Cells.Item(7, 7).Value = 2.45  # G7: 2.15 -> 2.45
$ws.Cells.Item(7, 8).Value = 3.3  # H7: 3.4 -> 3.3
$ws.Cells.Item(7, 9).Value = 2.75  # I7: 3.2 -> 2.75
$ws.Cells.Item(7, 10).Value = 3  # J7: 2.75 -> 3
$ws.Cells.Item(7, 11).Value = 2.25  # K7: 2.3 -> 2.25
$ws.Cells.Item(7, 12).Value = 3.25  # L7: 3.6 -> 3.25
$ws.Cells.Item(7, 15).Value = 1.22  # O7: 1.2 -> 1.22
$ws.Cells.Item(7, 16).Value = 4  # P7: 4.33 -> 4
$ws.Cells.Item(7, 17).Value = 1.75  # Q7: 1.7 -> 1.75
$ws.Cells.Item(7, 18).Value = 2.05  # R7: 2.1 -> 2.05
$ws.Cells.Item(7, 19).Value = 1.3  # S7: 1.33 -> 1.3
$ws.Cells.Item(7, 24).Value = 13  # X7: 12 -> 13
$ws.Cells.Item(7, 25).Value = 10  # Y7: 9 -> 10
$ws.Cells.Item(7, 26).Value = 23  # Z7: 21 -> 23
$ws.Cells.Item(7, 27).Value = 19  # AA7: 17 -> 19
$ws.Cells.Item(7, 29).Value = 12  # AC7: 13 -> 12
$ws.Cells.Item(7, 33).Value = 151  # AG7: 126 -> 151
$ws.Cells.Item(7, 34).Value = 11  # AH7: 13 -> 11
$ws.Cells.Item(7, 35).Value = 15  # AI7: 17 -> 15
$ws.Cells.Item(7, 36).Value = 11  # AJ7: 12 -> 11
$ws.Cells.Item(7, 37).Value = 29  # AK7: 34 -> 29
$ws.Cells.Item(7, 38).Value = 21  # AL7: 23 -> 21
$ws.Cells.Item(7, 39).Value = 26  # AM7: 29 -> 26
$ws.Cells.Item(7, 40).Value = 4.75  # AN7: 4.33 -> 4.75
$ws.Cells.Item(7, 41).Value = 13  # AO7: 11 -> 13
$ws.Cells.Item(7, 42).Value = 21  # AP7: 19 -> 21
$ws.Cells.Item(7, 45).Value = 126  # AS7: 101 -> 126
$ws.Cells.Item(7, 50).Value = 5  # AX7: 5.5 -> 5
$ws.Cells.Item(7, 51).Value = 15  # AY7: 17 -> 15
$ws.Cells.Item(7, 53).Value = 41  # BA7: 51 -> 41
$ws.Cells.Item(7, 54).Value = 51  # BB7: 67 -> 51

# Row 8
$ws.Cells.Item(8, 19).Value = 1.37  # S8: 1.4 -> 1.37

# Row 9
$ws.Cells.Item(9, 7).Value = 5.5  # G9: 6 -> 5.5
$ws.Cells.Item(9, 8).Value = 4  # H9: 4.1 -> 4
$ws.Cells.Item(9, 9).Value = 1.57  # I9: 1.53 -> 1.57
$ws.Cells.Item(9, 13).Value = 1.05  # M9: 1.04 -> 1.05
$ws.Cells.Item(9, 14).Value = 11  # N9: 13 -> 11
$ws.Cells.Item(9, 17).Value = 1.88  # Q9: 1.85 -> 1.88
$ws.Cells.Item(9, 18).Value = 1.98  # R9: 2 -> 1.98
$ws.Cells.Item(9, 19).Value = 1.33  # S9: 1.36 -> 1.33
$ws.Cells.Item(9, 24).Value = 29  # X9: 34 -> 29
$ws.Cells.Item(9, 25).Value = 17  # Y9: 19 -> 17
$ws.Cells.Item(9, 26).Value = 51  # Z9: 67 -> 51
$ws.Cells.Item(9, 28).Value = 41  # AB9: 51 -> 41
$ws.Cells.Item(9, 31).Value = 17  # AE9: 19 -> 17
$ws.Cells.Item(9, 33).Value = 301  # AG9: 351 -> 301
$ws.Cells.Item(9, 35).Value = 7.5  # AI9: 7 -> 7.5
$ws.Cells.Item(9, 39).Value = 26  # AM9: 29 -> 26
$ws.Cells.Item(9, 40).Value = 7  # AN9: 7.5 -> 7
$ws.Cells.Item(9, 41).Value = 29  # AO9: 34 -> 29
$ws.Cells.Item(9, 42).Value = 34  # AP9: 41 -> 34
$ws.Cells.Item(9, 43).Value = 101  # AQ9: 126 -> 101
$ws.Cells.Item(9, 44).Value = 126  # AR9: 151 -> 126
$ws.Cells.Item(9, 45).Value = 251  # AS9: 301 -> 251
$ws.Cells.Item(9, 51).Value = 8  # AY9: 7.5 -> 8

# Row 10
$ws.Cells.Item(10, 7).Value = 1.75  # G10: 1.85 -> 1.75
$ws.Cells.Item(10, 8).Value = 3.5  # H10: 3.4 -> 3.5
$ws.Cells.Item(10, 9).Value = 4.5  # I10: 4.1 -> 4.5
$ws.Cells.Item(10, 10).Value = 2.5  # J10: 2.6 -> 2.5
$ws.Cells.Item(10, 12).Value = 5  # L10: 4.75 -> 5
$ws.Cells.Item(10, 19).Value = 1.41  # S10: 1.44 -> 1.41
$ws.Cells.Item(10, 20).Value = 2.62  # T10: 2.63 -> 2.62
$ws.Cells.Item(10, 21).Value = 2  # U10: 1.91 -> 2
$ws.Cells.Item(10, 22).Value = 1.73  # V10: 1.8 -> 1.73
$ws.Cells.Item(10, 24).Value = 8  # X10: 8.5 -> 8
$ws.Cells.Item(10, 26).Value = 13  # Z10: 15 -> 13
$ws.Cells.Item(10, 27).Value = 15  # AA10: 17 -> 15
$ws.Cells.Item(10, 30).Value = 7  # AD10: 6.5 -> 7
$ws.Cells.Item(10, 37).Value = 51  # AK10: 41 -> 51
$ws.Cells.Item(10, 38).Value = 41  # AL10: 34 -> 41
$ws.Cells.Item(10, 41).Value = 9.5  # AO10: 10 -> 9.5
$ws.Cells.Item(10, 42).Value = 21  # AP10: 23 -> 21
$ws.Cells.Item(10, 48).Value = 67  # AV10: 51 -> 67
$ws.Cells.Item(10, 51).Value = 26  # AY10: 23 -> 26
$ws.Cells.Item(10, 54).Value = 126  # BB10: 101 -> 126
$ws.Cells.Item(10, 55).Value = 301  # BC10: 251 -> 301

# Row 14
$ws.Cells.Item(14, 9).Value = 2.6  # I14: 2.63 -> 2.6
$ws.Cells.Item(14, 10).Value = 3.4  # J14: 3.25 -> 3.4
$ws.Cells.Item(14, 50).Value = 4.5  # AX14: 4.75 -> 4.5

# Row 15
$ws.Cells.Item(15, 7).Value = 2.9  # G15: 2.75 -> 2.9
$ws.Cells.Item(15, 8).Value = 3.25  # H15: 3.3 -> 3.25
$ws.Cells.Item(15, 9).Value = 2.38  # I15: 2.5 -> 2.38
$ws.Cells.Item(15, 10).Value = 3.5  # J15: 3.4 -> 3.5
$ws.Cells.Item(15, 13).Value = 1.06  # M15: 1.05 -> 1.06
$ws.Cells.Item(15, 14).Value = 10  # N15: 11 -> 10
$ws.Cells.Item(15, 21).Value = 1.73  # U15: 1.67 -> 1.73
$ws.Cells.Item(15, 22).Value = 2  # V15: 2.1 -> 2
$ws.Cells.Item(15, 25).Value = 11  # Y15: 10 -> 11
$ws.Cells.Item(15, 27).Value = 23  # AA15: 21 -> 23
$ws.Cells.Item(15, 29).Value = 10  # AC15: 11 -> 10
$ws.Cells.Item(15, 34).Value = 8.5  # AH15: 9 -> 8.5
$ws.Cells.Item(15, 35).Value = 12  # AI15: 13 -> 12
$ws.Cells.Item(15, 36).Value = 9.5  # AJ15: 10 -> 9.5
$ws.Cells.Item(15, 40).Value = 5  # AN15: 4.75 -> 5

Write-Output "Applied 162 cell updates"